# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G (K = strikeouts) for rows 2-19 with the recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 3
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
